$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Clear out the old "per-day" data blocks that are being removed
#    (rows 2-3 B:F/H, and row 8 B/D:F). Styles (borders/fills) stay.
# ---------------------------------------------------------------
$ws.Range("B2:F2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("B3:F3").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("D8:F8").ClearContents()

# Old helper formulas in I3/J3 (I2/18+0.4 style) are no longer used
$ws.Range("I3:J3").ClearContents()

# ---------------------------------------------------------------
# 2. Row 2: I2 becomes a literal, J2:M2 become new SUM() formulas
# ---------------------------------------------------------------
$ws.Range("I2").Value = 30
$ws.Range("J2").Formula = "=SUM(C3:C6)"
$ws.Range("K2").Formula = "=SUM(D3:D6)"
$ws.Range("L2").Formula = "=SUM(E3:E6)"
$ws.Range("M2").Formula = "=SUM(F3:F6)"
$ws.Range("M2").NumberFormat = $ws.Range("L2").NumberFormat

# ---------------------------------------------------------------
# 3. Row 3: K3:M3 become new SUM() formulas
# ---------------------------------------------------------------
$ws.Range("K3").Formula = "=SUM(D9:D11)"
$ws.Range("L3").Formula = "=SUM(E9:E11)"
$ws.Range("M3").Formula = "=SUM(F9:F11)"
$ws.Range("M3").NumberFormat = $ws.Range("L3").NumberFormat

# ---------------------------------------------------------------
# 4. Row 5: new subtotal formulas J5:M5 = SUM() of rows 2:4
# ---------------------------------------------------------------
$ws.Range("J5").Formula = "=SUM(J2:J4)"
$ws.Range("K5").Formula = "=SUM(K2:K4)"
$ws.Range("L5").Formula = "=SUM(L2:L4)"
$ws.Range("M5").Formula = "=SUM(M2:M4)"
$ws.Range("J5:M5").NumberFormat = $ws.Range("I2").NumberFormat

# ---------------------------------------------------------------
# 5. Row 6: new literal "used" values J6:M6 (plain, unstyled)
# ---------------------------------------------------------------
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 101
$ws.Range("L6").Value = 17
$ws.Range("M6").Value = 7

# ---------------------------------------------------------------
# 6. Row 7: new "remaining" formulas J7:M7 = row5 - row6
# ---------------------------------------------------------------
$ws.Range("J7").Formula = "=J5-J6"
$ws.Range("K7").Formula = "=K5-K6"
$ws.Range("L7").Formula = "=L5-L6"
$ws.Range("M7").Formula = "=M5-M6"
$ws.Range("J7:M7").NumberFormat = $ws.Range("I2").NumberFormat

# ---------------------------------------------------------------
# 7. Row 8: new per-unit formulas J8:M8 = row7 / 19
# ---------------------------------------------------------------
$ws.Range("J8").Formula = "=J7/19"
$ws.Range("K8").Formula = "=K7/19"
$ws.Range("L8").Formula = "=L7/19"
$ws.Range("M8").Formula = "=M7/19"
$ws.Range("J8:M8").NumberFormat = $ws.Range("I3").NumberFormat

# ---------------------------------------------------------------
# 8. Row 9: I9 literal divisor (style copied from B2 minus border),
#    J9:M9 = row8 / $I9
# ---------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("I9").Borders.LineStyle = -4142
$ws.Range("I9").Value = 10

$ws.Range("J9").Formula = "=J8/`$I9"
$ws.Range("K9").Formula = "=K8/`$I9"
$ws.Range("L9").Formula = "=L8/`$I9"
$ws.Range("M9").Formula = "=M8/`$I9"
$ws.Range("J9:M9").NumberFormat = $ws.Range("I3").NumberFormat

# ---------------------------------------------------------------
# 9. Sheet view: selection moves to K11, no frozen/topLeft scroll cell
# ---------------------------------------------------------------
$ws.Range("K11").Select()
